$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------------
# Column F ("Include in Occupancy Calculation") flips from "Yes" to "No"
# for a handful of rows (Greenwood Village, Tampa, Santiago, Sao Paulo).
$ws.Range("F16").Value = "No"
$ws.Range("F38").Value = "No"
$ws.Range("F47").Value = "No"
$ws.Range("F48").Value = "No"

# Melbourne desk count (column C, row 44) bumps from 30 to 32.
$ws.Range("C44").Value = 32

# --- View/selection state ---------------------------------------------------
# Scroll the sheet and move the active selection, matching the saved
# workbook view (best-effort: zoom + selection are supported; topLeftCell
# scroll position is not exposed for persistence in this runtime).
$win = $excel.ActiveWindow
$ws.Range("D46").Select() | Out-Null
$win.ScrollRow = 22
$win.ScrollColumn = 1
$win.Zoom = 87
